$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text (never auto-converted to a number/date),
# without altering the cells NumberFormat/Style (mirrors how the source data keeps
# numeric-looking price strings as plain text/inline strings).
function Set-TextValue($cellRef, [string]$val) {
    $escaped = $val.Replace('"', '""')
    $ws.Range($cellRef).Formula = '="' + $escaped + '"'
    $ws.Range($cellRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)  # xlPasteValues
}

# Row 2
Set-TextValue "D2" '97.554.03'
$ws.Range("E2").Value = '  -1.75%  '

# Row 3
Set-TextValue "D3" '3.411.72'
$ws.Range("E3").Value = '  +4.02%  '

# Row 4
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
Set-TextValue "D5" '254.93'
$ws.Range("E5").Value = '  +0.60%  '

# Row 6
Set-TextValue "D6" '653.37'

# Row 7
$ws.Range("E7").Value = '  +3.78%  '

# Row 8
$ws.Range("E8").Value = '  +6.15%  '

# Row 9
$ws.Range("E9").Value = '  +9.98%  '

# Row 10
$ws.Range("E10").Value = '  -0.03%  '

# Row 11
Set-TextValue "D11" '3.407.10'
$ws.Range("E11").Value = '  +3.92%  '

# Row 12
$ws.Range("E12").Value = '  +6.21%  '

# Row 13
Set-TextValue "D13" '41.80'
$ws.Range("E13").Value = '  +6.50%  '

# Row 14
$ws.Range("E14").Value = '  +14.81%  '

# Row 15
$ws.Range("E15").Value = '  +5.34%  '

# Row 16
Set-TextValue "D16" '97.229.76'

# Row 17
Set-TextValue "D17" '4.045.47'
$ws.Range("E17").Value = '  +4.56%  '

# Row 18
$ws.Range("E18").Value = '  +37.21%  '

# Row 19
Set-TextValue "D19" '3.413.39'
$ws.Range("E19").Value = '  +4.44%  '

# Row 20
Set-TextValue "D20" '17.45'
$ws.Range("E20").Value = '  +14.50%  '

# Row 21
Set-TextValue "D21" '0.510'
$ws.Range("E21").Value = '  +58.82%  '

# Row 22
$ws.Range("E22").Value = '  +17.42%  '

# Row 23
$ws.Range("E23").Value = '  +0.85%  '

# Row 24
Set-TextValue "D24" '504.82'
$ws.Range("E24").Value = '  +3.74%  '

# Row 25
$ws.Range("E25").Value = '  +2.49%  '

# Row 26
$ws.Range("E26").Value = '  +9.19%  '

# Row 27
Set-TextValue "D27" '98.83'
$ws.Range("E27").Value = '  +11.07%  '

# Row 28
Set-TextValue "D28" '12.71'
$ws.Range("E28").Value = '  +6.30%  '

# Row 29
Set-TextValue "D29" '3.599.19'
$ws.Range("E29").Value = '  +5.12%  '

# Row 30
Set-TextValue "D30" '0.156'
$ws.Range("E30").Value = '  +14.99%  '

# Row 31
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue "D31" '11.40'
$ws.Range("E31").Value = '  +10.74%  '

# Row 32
$ws.Range("B32").Value = 'Cronos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue "D32" '0.198'
$ws.Range("E32").Value = '  +4.88%  '

# Row 33
Set-TextValue "D33" '0.998'
$ws.Range("E33").Value = '  -0.19%  '

# Row 34
Set-TextValue "D34" '1.00'
$ws.Range("E34").Value = '  +0.16%  '

# Row 35
$ws.Range("E35").Value = '  +21.63%  '

# Row 36
Set-TextValue "D36" '29.83'
$ws.Range("E36").Value = '  +7.08%  '

# Row 37
Set-TextValue "D37" '2.26'
$ws.Range("E37").Value = '  +17.27%  '

# Row 38
Set-TextValue "D38" '7.74'
$ws.Range("E38").Value = '  +7.90%  '

# Row 39
$ws.Range("E39").Value = '  +2.83%  '

# Row 40
$ws.Range("E40").Value = '  +15.45%  '

# Row 41
Set-TextValue "D41" '513.10'
$ws.Range("E41").Value = '  +5.75%  '

# Row 42
$ws.Range("E42").Value = '  -0.42%  '

# Row 43
Set-TextValue "D43" '0.856'
$ws.Range("E43").Value = '  +11.27%  '

# Row 44
Set-TextValue "D44" '0.0421'
$ws.Range("E44").Value = '  +25.70%  '

# Row 45
Set-TextValue "D45" '3.67'
$ws.Range("E45").Value = '  +0.02%  '

# Row 46
$ws.Range("B46").Value = 'Filecoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-TextValue "D46" '5.46'
$ws.Range("E46").Value = '  +16.28%  '

# Row 47
$ws.Range("B47").Value = 'dogwifhat'
$ws.Range("C47").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue "D47" '3.27'
$ws.Range("E47").Value = '  +6.55%  '

# Row 48
$ws.Range("E48").Value = '  +0.00%  '

# Row 49
Set-TextValue "D49" '8.16'
$ws.Range("E49").Value = '  +12.28%  '

# Row 50
$ws.Range("E50").Value = '  +16.52%  '

# Row 51
$ws.Range("E51").Value = '  +6.93%  '
